# Fix imported devices template: remove the obsolete "索引" (Index) column
# and correct the example row's auto-subscribe value from "0/1" to "0".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("索引") is no longer needed; delete the whole column so that
# everything to its right (所属网关, 设备编号, ...) shifts one column left.
$ws.Columns.Item(6).Delete()

# After the shift, the example row's auto-subscribe example value
# (originally in column S, now column R) should read "0" instead of "0/1".
$ws.Range("R3").Value = "0"
